# Flare "Monsters" data import — enable can_cast / can_use_artifacts
# (columns M and N) for every monster that didn't already have them set,
# and leave the sheet scrolled/selected on the newly-populated column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monsters")

# Rows 2-63 and 73-106 are missing the "can_cast" (M) / "can_use_artifacts"
# (N) flags; every other populated row in the sheet already has them set
# to 1, so bring these in line.
$rowRanges = @(@(2, 63), @(73, 106))

foreach ($rowRange in $rowRanges) {
    $startRow = $rowRange[0]
    $endRow = $rowRange[1]
    for ($r = $startRow; $r -le $endRow; $r++) {
        $ws.Cells.Item($r, 13).Value = 1
        $ws.Cells.Item($r, 14).Value = 1
    }
}

# Leave the sheet scrolled back up with the new "can_use_artifacts" column
# selected, rather than where the previous editor had left it (N70/Z103).
$ws.Activate()
$ws.Range("N2:N111").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
